# Commit: "Mon, Mar 30, 2020  3:05:31 AM"
#
# The underlying diff changes the table style used by the three summary
# tables (balance-sheet style tables on the "Fixed Assets" slides) from the
# deck's custom table style {EAB39EEA-4CAF-4D4D-AD32-453E067FD839} to the
# built-in gallery style {DD68FA76-7987-4B9C-92E8-E99884183B14}.
#
# In the PowerPoint UI this is exactly what happens when you select a table
# and pick a different style from Table Tools > Design > Table Styles -
# which is exposed on the object model as Table.ApplyStyle(StyleID).
# (Table.Style is read-only / cannot be assigned directly.)

$p = $ppt.ActivePresentation

$oldStyleId = "{EAB39EEA-4CAF-4D4D-AD32-453E067FD839}"
$newStyleId = "{DD68FA76-7987-4B9C-92E8-E99884183B14}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shp = $slide.Shapes.Item($shi)
        if ($shp.HasTable) {
            $tbl = $shp.Table
            if ($tbl.Style -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}
